# modificando ejemplos TEI 0.2.1
# Fill the "Definition" column (D) on the "Concepts" sheet with the same
# values as the "Display" column (C) for each concept row (2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 11; $row++) {
    $display = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value2 = $display
}
